# Calibrate Georgia for total population
# Update the "susceptible_fully" starting population constant from
# 4,000,000 to 3,700,000 on the "constants" sheet, and leave the
# active selection on the edited cell (B4), matching the author's
# edit captured in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

$ws.Activate()
$ws.Range("B4").Value = 3700000
$ws.Range("B4").Select()
